$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Swap the contents of column C (codeforiati:group-code) and column D (codeforiati:group-name)
# for every row, including the header row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
